$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "XPT1" sheet at the very end of the workbook first, so that
#    it receives sheetId 4 (the workbook's next free id at that point).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newXpt1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newXpt1.Name = "XPT1"

# ---------------------------------------------------------------------------
# 2. Add the new "Spirit of QLD" sheet right after "Stations" (it becomes the
#    second tab). Added after XPT1 so it receives sheetId 5.
# ---------------------------------------------------------------------------
$stationsTmp = $wb.Worksheets.Item("Stations")
$newQld = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $stationsTmp)
$newQld.Name = "Spirit of QLD"

# ---------------------------------------------------------------------------
# Re-fetch every sheet handle by name now that the tab order/positions are
# final -- object references captured before/around an Add() call track the
# tab *index*, not the sheet identity, so they can silently go stale once
# later inserts shift everything after them.
# ---------------------------------------------------------------------------
$stations = $wb.Worksheets.Item("Stations")
$qld = $wb.Worksheets.Item("Spirit of QLD")
$overland = $wb.Worksheets.Item("Overland Line")
$trainServices = $wb.Worksheets.Item("TrainServices")
$xpt1 = $wb.Worksheets.Item("XPT1")

# ---------------------------------------------------------------------------
# 3. Populate "Spirit of QLD" (Brisbane -> Cairns order list).
# ---------------------------------------------------------------------------
$qldStations = @(
    "Brisbane",
    "Caboolture",
    "Landsborough",
    "Nambour",
    "Cooroy",
    "Gympie North",
    "Maryborough West",
    "Howard",
    "Bundaberg",
    "Miriam Vale",
    "Gladstone",
    "Mount Larcom",
    "Rockhampton",
    "St Lawrence",
    "Carmila",
    "Sarina",
    "Mackay",
    "Proserpine",
    "Bowen",
    "Home Hill",
    "Ayr",
    "Giru",
    "Townsville",
    "Ingham",
    "Cardwell",
    "Tully",
    "Innisfail",
    "Babinda",
    "Gordonvale",
    "Cairns"
)

$qld.Range("A1").Value = "Station_name"
$qld.Range("B1").Value = "Order"
for ($i = 0; $i -lt $qldStations.Count; $i++) {
    $row = $i + 2
    $qld.Cells.Item($row, 1).Value = $qldStations[$i]
    $qld.Cells.Item($row, 2).Value = $i + 1
}

$qld.Columns.Item(1).ColumnWidth = 16.77734375

# ---------------------------------------------------------------------------
# 4. Populate "XPT1" (Sydney -> Melbourne order list).
# ---------------------------------------------------------------------------
$xptStations = @(
    "Central Station",
    "Campbelltown Station",
    "Moss Vale Station",
    "Goulburn Station",
    "Gunning Station",
    "Yass Junction station",
    "Harden Station",
    "Cootamundra Station",
    "Junee Station",
    "Wagga Wagga Station",
    "The Rock Station",
    "Henty Station",
    "Culcairn Station",
    "Albury Station",
    "Wangaratta Station",
    "Benella Station",
    "Seymour Station",
    "Broadmeadows Station",
    "Melbourne Southern Cross"
)

$xpt1.Range("A1").Value = "Station_name"
$xpt1.Range("B1").Value = "Order"
for ($i = 0; $i -lt $xptStations.Count; $i++) {
    $row = $i + 2
    $xpt1.Cells.Item($row, 1).Value = $xptStations[$i]
    $xpt1.Cells.Item($row, 2).Value = $i + 1
}

$xpt1.Columns.Item(1).ColumnWidth = 25.44140625
[void]$xpt1.Activate()
[void]$xpt1.Range("A20").Select()

# ---------------------------------------------------------------------------
# 5. Minor view-state tweaks on the pre-existing sheets.
# ---------------------------------------------------------------------------
[void]$overland.Activate()
[void]$overland.Range("D17").Select()

[void]$trainServices.Activate()
[void]$trainServices.Range("H3").Select()

[void]$stations.Activate()
[void]$stations.Range("B51:B80").Select()

# Leave "Spirit of QLD" as the active/selected tab (matches activeTab="1",
# i.e. the second tab - "Spirit of QLD" - in the final sheet order).
[void]$qld.Activate()
[void]$qld.Range("G19").Select()
